$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove trailing rows 24 and 25 (their content is merged into rows 22/23 below)
$ws.Rows("24:25").Delete()

# Row 10
$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "3577649 - Carlos Angelo Nunes"
$ws.Range("C10").Value = "3577649 - Carlos Angelo Nunes"
$ws.Range("A10").RowHeight = 60

# Row 11
$ws.Range("A11").Value = "Objectives:"
$ws.Range("B11").ClearContents()
$ws.Range("C11").ClearContents()
$ws.Range("A11").RowHeight = 60

# Row 12
$ws.Range("A12").Value = "Docentes responsáveis:"
$ws.Range("B12").ClearContents()
$ws.Range("C12").ClearContents()
$ws.Rows("12").AutoFit()

# Row 13
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "01/01/2022"
$ws.Range("C13").Value = "01/01/2022"
$ws.Range("A13").RowHeight = 60

# Row 14
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Range("A14").RowHeight = 60

# Row 15
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "3577649 - Carlos Angelo Nunes"
$ws.Range("C15").Value = "3577649 - Carlos Angelo Nunes"
$ws.Range("A15").RowHeight = 120

# Row 16
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").ClearContents()
$ws.Range("C16").ClearContents()
$ws.Range("A16").RowHeight = 120

# Row 17
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Rows("17").AutoFit()

# Row 18
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "1922320 - Sebastiao Ribeiro"
$ws.Range("C18").Value = "1922320 - Sebastiao Ribeiro"
$ws.Range("A18").RowHeight = 60

# Row 19
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "O curso será ministrado na forma de aulas expositivas."
$ws.Range("C19").Value = "O curso será ministrado na forma de aulas expositivas."
$ws.Range("A19").RowHeight = 60

# Row 20
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Serão aplicadas duas avaliações escritas (P1 e P2) que comporão a nota final (NF). A nota final será calculada através da expressão: NF=(P1+P2)/2. Poderão ser solicitados trabalhos aos alunos e que irão também compor parte da nota P2."
$ws.Range("C20").Value = "Serão aplicadas duas avaliações escritas (P1 e P2) que comporão a nota final (NF). A nota final será calculada através da expressão: NF=(P1+P2)/2. Poderão ser solicitados trabalhos aos alunos e que irão também compor parte da nota P2."
$ws.Range("A20").RowHeight = 60

# Row 21
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Para os alunos que obtiverem 3,0"
$ws.Range("C21").Value = "Para os alunos que obtiverem 3,0"
$ws.Range("A21").RowHeight = 120

# Row 22
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Rows("22").AutoFit()

# Row 23
$ws.Range("A23").ClearContents()
$ws.Range("B23").Value = "LOM3015 -  Termodinâmica de Materiais  (Requisito fraco)`n"
$ws.Range("C23").Value = "LOM3015 -  Termodinâmica de Materiais  (Requisito fraco)`n"
$ws.Range("A23").RowHeight = 30
